# ----------------------------------------------------------------------------
# LOB1223.docx edit script
#
# 1) "Curso (semestre ideal): EA (7)" -> "Curso (semestre ideal): EA (9)"
# 2) The Requisitos bullet paragraph's three "(Requisito)" lines are replaced
#    by the expanded list of 28 "(Requisito fraco)" lines (same ListBullet
#    paragraph, each entry kept as its own run terminated by a manual
#    line break, matching the original run layout).
# ----------------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Curso (semestre ideal): EA (7) -> EA (9) ---------------------------
$d.Content.Find.Execute("Curso (semestre ideal): EA (7)", $false, $false, $false, $false, $false, `
    $true, 1, $false, "Curso (semestre ideal): EA (9)", 2) | Out-Null

# --- 2) Rebuild the Requisitos bullet list ----------------------------------
# Locate the ListBullet paragraph that still holds the old "LOB1210" entry
# (falls back to the document's last paragraph if, for any reason, the text
# search below does not hit).
$reqPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*LOB1210*") {
        $reqPara = $para
    }
}
if ($null -eq $reqPara) {
    $reqPara = $d.Paragraphs.Last
}
$reqRange = $reqPara.Range

$flatOpc = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>LOB1003 -  Cálculo I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1004 -  Cálculo II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1006 -  Cálculo IV  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1011 -  Eletricidade Aplicada  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1012 -  Estatística  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1018 -  Física I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1019 -  Física II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1021 -  Física IV  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1024 -  Mecânica  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1036 -  Geometria Analítica  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1037 -  Àlgebra Linear  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1038 -  Física Experimental I  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1039 -  Física Experimental III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1041 -  Física Experimental II  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1042 -  Física Experimental IV  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1052 -  Cálculo III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1053 -  Física III  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1232 -  Licenciamento Ambiental  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOB1257 -  Sistema de Abastecimento e Tratamento de Água  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOM3081 -  Introdução à Mecânica dos Sólidos  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4095 -  Química Geral Experimental  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4097 -  Fundamentos de Química para Engenharia I (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4233 -  Gestão de Negócios  (Requisito fraco)</w:t><w:br/></w:r><w:r><w:t>LOQ4247 -  Desenho Assistido por Computador  (Requisito fraco)</w:t><w:br/></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

# InsertXML REPLACES $reqRange's contents with the supplied WordprocessingML.
# Because the payload is a full <w:p>, this lands as a brand new paragraph
# and pushes the original (now empty) ListBullet paragraph mark into a
# trailing empty paragraph immediately afterward.
$reqRange.InsertXML($flatOpc) | Out-Null

# Merge that trailing empty paragraph back into the one we just inserted by
# deleting the paragraph mark that now separates them; the surviving mark
# (and therefore its ListBullet style) is the original paragraph's.
$newParaIndex = $d.Paragraphs.Count - 1
$newPara = $d.Paragraphs($newParaIndex)
$markStart = $newPara.Range.End - 1
$d.Range($markStart, $markStart + 1).Delete() | Out-Null
